# Update "想去人数" (F column) figures across sheets to the newly scraped
# counts. Same underlying events appear on multiple tabs (展览 / 演出 /
# 全部类型), so every occurrence of a given row is updated together.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 863
$ws1.Range("F8").Value  = 50
$ws1.Range("F12").Value = 908
$ws1.Range("F16").Value = 413
$ws1.Range("F17").Value = 6631
$ws1.Range("F21").Value = 7580
$ws1.Range("F26").Value = 1799
$ws1.Range("F27").Value = 894
$ws1.Range("F29").Value = 130
$ws1.Range("F31").Value = 71
$ws1.Range("F32").Value = 220
$ws1.Range("F33").Value = 197
$ws1.Range("F34").Value = 1675
$ws1.Range("F40").Value = 1762
$ws1.Range("F41").Value = 2135

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 80

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value  = 863
$ws4.Range("F10").Value = 50
$ws4.Range("F14").Value = 909
$ws4.Range("F19").Value = 413
$ws4.Range("F20").Value = 6631
$ws4.Range("F24").Value = 7580
$ws4.Range("F29").Value = 1799
$ws4.Range("F30").Value = 894
$ws4.Range("F32").Value = 130
$ws4.Range("F34").Value = 71
$ws4.Range("F36").Value = 220
$ws4.Range("F37").Value = 197
$ws4.Range("F38").Value = 1675
$ws4.Range("F45").Value = 1762
$ws4.Range("F47").Value = 2135
$ws4.Range("F49").Value = 80
